$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.857.86"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.838.20"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.32"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4722"
$ws.Range("E7").Value = "  +3.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3654"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07145"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9181"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.965.85"
$ws.Range("E11").Value = "  +8.90%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.52"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07662"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.281"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.392"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.97"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008620"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.882.05"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.46"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.005"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.920"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.72"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.20"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.003"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.14"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.864"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08822"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.211"
$ws.Range("E31").Value = "  +3.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.174"
$ws.Range("E32").Value = "  +5.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7431"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.469"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.741"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.088"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05214"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.958"
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5187"
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.958"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1511"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.139"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.40"
$ws.Range("E44").Value = "  +4.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4695"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.006"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.37"
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.591"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.83"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06029"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8853"
$ws.Range("E51").Value = "  +4.73%  "
